# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Mon Feb 12 21:39:03 UTC 2024 with GitHub Actions".
# For every changed row we update the Price (column D) and/or the
# Volume(1h) (column E) cell; row 29/30 additionally swap Kaspa and Cosmos.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bitcoin (row 2)
$ws.Range("D2").Value = '50.018.98'
$ws.Range("E2").Value = '  +4.18%  '

# Ethereum (row 3)
$ws.Range("D3").Value = '2.638.47'
$ws.Range("E3").Value = '  +5.57%  '

# BNB (row 5)
$ws.Range("D5").Value = '''327.69'
$ws.Range("E5").Value = '  +2.26%  '

# Solana (row 6)
$ws.Range("D6").Value = '''110.60'
$ws.Range("E6").Value = '  +3.01%  '

# XRP (row 7)
$ws.Range("D7").Value = '''0.533'
$ws.Range("E7").Value = '  +1.73%  '

# USDC (row 8)
$ws.Range("E8").Value = '  +0.03%  '

# Cardano (row 9)
$ws.Range("D9").Value = '''0.561'
$ws.Range("E9").Value = '  +4.41%  '

# Avalanche (row 10)
$ws.Range("D10").Value = '''40.77'
$ws.Range("E10").Value = '  +3.22%  '

# Chainlink (row 11)
$ws.Range("D11").Value = '''20.68'
$ws.Range("E11").Value = '  +2.09%  '

# Dogecoin (row 12)
$ws.Range("D12").Value = '''0.0822'
$ws.Range("E12").Value = '  +1.31%  '

# TRON (row 13)
$ws.Range("E13").Value = '  +0.93%  '

# Polkadot (row 14)
$ws.Range("E14").Value = '  +2.76%  '

# WrappedliquidstakedEther2.0 (row 15)
$ws.Range("D15").Value = '3.052.84'
$ws.Range("E15").Value = '  +5.63%  '

# WrappedEther (row 16)
$ws.Range("D16").Value = '2.625.60'
$ws.Range("E16").Value = '  +5.55%  '

# Polygon (row 17)
$ws.Range("E17").Value = '  +5.38%  '

# WrappedBTC (row 18)
$ws.Range("D18").Value = '49.962.82'
$ws.Range("E18").Value = '  +4.38%  '

# ImmutableX (row 19)
$ws.Range("E19").Value = '  +11.61%  '

# InternetComputer(DFINITY) (row 20)
$ws.Range("D20").Value = '''13.35'
$ws.Range("E20").Value = '  +3.35%  '

# Uniswap (row 21)
$ws.Range("E21").Value = '  +1.63%  '

# ShibaInu (row 22)
$ws.Range("D22").Value = '0.0₃0964'
$ws.Range("E22").Value = '  +2.68%  '

# Litecoin (row 23)
$ws.Range("D23").Value = '''73.02'
$ws.Range("E23").Value = '  +2.22%  '

# BitcoinCash (row 24)
$ws.Range("D24").Value = '''280.44'
$ws.Range("E24").Value = '  +0.89%  '

# PancakeSwap (row 25)
$ws.Range("E25").Value = '  +2.43%  '

# EthereumClassic (row 26)
$ws.Range("D26").Value = '''26.64'
$ws.Range("E26").Value = '  +4.12%  '

# Dai (row 27)
$ws.Range("E27").Value = '  -0.10%  '

# Toncoin (row 28)
$ws.Range("D28").Value = '''2.24'
$ws.Range("E28").Value = '  +1.67%  '

# Kaspa (row 29)
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").Value = '''9.95'
$ws.Range("E29").Value = '  +2.35%  '

# Cosmos (row 30)
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").Value = '''0.145'
$ws.Range("E30").Value = '  +3.97%  '

# InjectiveProtocol (row 31)
$ws.Range("E31").Value = '  +4.61%  '

# OKB (row 32)
$ws.Range("D32").Value = '''49.75'
$ws.Range("E32").Value = '  +0.69%  '

# Celestia (row 33)
$ws.Range("D33").Value = '''19.78'
$ws.Range("E33").Value = '  +1.54%  '

# Filecoin (row 34)
$ws.Range("E34").Value = '  +2.94%  '

# Hedera (row 36)
$ws.Range("E36").Value = '  +2.23%  '

# ARBITRUM (row 37)
$ws.Range("D37").Value = '''2.07'
$ws.Range("E37").Value = '  +6.69%  '

# RenderToken (row 38)
$ws.Range("E38").Value = '  +3.29%  '

# LidoDAOToken (row 39)
$ws.Range("D39").Value = '''3.11'
$ws.Range("E39").Value = '  +7.65%  '

# Stellar (row 40)
$ws.Range("E40").Value = '  +1.50%  '

# Monero (row 41)
$ws.Range("D41").Value = '''123.44'
$ws.Range("E41").Value = '  +2.37%  '

# EnergySwap (row 42)
$ws.Range("D42").Value = '''22.61'
$ws.Range("E42").Value = '  +5.97%  '

# WEMIXToken (row 43)
$ws.Range("E43").Value = '  +0.37%  '

# VeChain (row 44)
$ws.Range("E44").Value = '  +4.73%  '

# NEARProtocol (row 45)
$ws.Range("E45").Value = '  +6.31%  '

# Maker (row 46)
$ws.Range("D46").Value = '2.065.87'
$ws.Range("E46").Value = '  +2.96%  '

# ApeXProtocol (row 47)
$ws.Range("D47").Value = '''2.31'
$ws.Range("E47").Value = '  +15.64%  '

# Stacks (row 48)
$ws.Range("E48").Value = '  +8.78%  '

# FraxShare (row 49)
$ws.Range("E49").Value = '  +1.18%  '

# THORChain (row 50)
$ws.Range("D50").Value = '''5.39'
$ws.Range("E50").Value = '  +4.48%  '

# BitcoinSV (row 51)
$ws.Range("D51").Value = '''81.80'
$ws.Range("E51").Value = '  +2.16%  '
